# Edit script: apply documentation.docx changes
$d = $word.ActiveDocument

# --- Step 1: Remove the stray "_GoBack" bookmark that currently sits at the
#     end of the "Routes" paragraph (it will be recreated further up, inside
#     the Overview paragraph, in step 3). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: Title paragraph - merge "Express" / " " / "Api" (proofed) /
#     " Server Template Documentation" runs into a single plain run. ---
$titleRng = $d.Content
$titleRng.Find.Execute("Express Api Server Template Documentation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$titleStart = $titleRng.Start
$titleEnd = $titleRng.End
$titleTarget = $d.Range($titleStart, $titleEnd)
$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Express Api Server Template Documentation</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$titleTarget.InsertXML($titleXml)

# --- Step 3: Overview paragraph - capitalize "WaBiSQue" correctly, drop the
#     spell-check proofing marks, and extend the sentence with the new
#     description (including the relocated "_GoBack" bookmark). ---
$overviewRng = $d.Content
$overviewRng.Find.Execute("The Wabisque Express Api Server Template (WEAST)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$overviewStart = $overviewRng.Start
$overviewEnd = $overviewRng.End
$overviewTarget = $d.Range($overviewStart, $overviewEnd)
$overviewXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">The </w:t></w:r><w:r><w:t>WaBiSQue</w:t></w:r><w:r><w:t xml:space="preserve"> Express Api Server Template (WEAST)</w:t></w:r><w:r><w:t xml:space="preserve"> is a simple template</w:t></w:r><w:r><w:t xml:space="preserve"> project</w:t></w:r><w:r><w:t xml:space="preserve"> for creating api servers with</w:t></w:r><w:r><w:t xml:space="preserve"> express and node. The server can be started by running the ‘index.js’ file located at the root directory of the project – all configurations are also mad</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>e in the ‘index.js’ file.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$overviewTarget.InsertXML($overviewXml)
